# WatchMe/KickJr PPG-Module BOM update
# - Replaces the "example" BOM rows (3-6) with the real PPG module parts
# - Drops the two extra rows (7-8) that existed in the old sheet
# - Removes the now-unused "this is just an example" comments on B3:B8
#   (keeps the still-relevant header comments on B2:E2 and J2)
# - Updates the title cell and the two SUM() totals to match the new extent

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- remove the stale example comments on column B (rows 3-8) ----
foreach ($ref in @("B3", "B4", "B5", "B6", "B7", "B8")) {
    $cmt = $ws.Range($ref).Comment
    if ($cmt -ne $null) {
        $cmt.Delete()
    }
}

# ---- drop rows 7 and 8 entirely (sheet shrinks from A1:O8 to A1:O6) ----
$ws.Rows("7:8").Delete()

# ---- title / totals row ----
$ws.Range("A1").Value = "Kick Jr PPG Module Rev A / Purdue Milestones"
$ws.Range("M1").Formula = "=SUM(M3:M6)"
$ws.Range("O1").Formula = "=SUM(O3:O6)"

# ---- row 3: Photodiode ----
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "Photodiode"
$ws.Range("C3").Value = "PD"
$ws.Range("D3").Value = "Photodiode 940nm 100ns 130° 2-SMD, Gull Wing, PHOTODIODE PIN HI SPEED HI SENS"
$ws.Range("E3").Value = "SMD"
$ws.Range("F3").Value = "2-SMD, Gull Wing"
$ws.Range("G3").Value = "Vishay Semiconductor Opto Division"
$ws.Range("H3").Value = "VBPW34S"
$ws.Range("I3").Value = "Digi-Key"
$ws.Range("J3").Value = "'751-1500-1-ND"
$ws.Range("J3").WrapText = $true
$ws.Range("K3").Value = 1.19
$ws.Range("L3").Value = 1
$ws.Range("M3").Formula = "=K3*L3"
$ws.Range("N3").Value = 0.566
$ws.Range("O3").Formula = "=N3*L3"

# ---- row 4: Green LED ----
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "Green LED"
$ws.Range("C4").Value = "LED1"
$ws.Range("D4").Value = "Green 525nm LED Indication - Discrete 3.1V 0805 (2012 Metric)"
$ws.Range("E4").Value = "SMD"
$ws.Range("F4").Value = "0805"
$ws.Range("G4").Value = "QT Brightek (QTB)"
$ws.Range("H4").Value = "QBLP631-IG"
$ws.Range("I4").Value = "Digi-Key"
$ws.Range("J4").Value = "1516-1084-1-ND"
$ws.Range("K4").Value = 0.39
$ws.Range("L4").Value = 1
$ws.Range("M4").Formula = "=K4*L4"
$ws.Range("N4").Value = 0.1617
$ws.Range("O4").Formula = "=N4*L4"

# ---- row 5: 2.2uF capacitor ----
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "2.2µF capacitor"
$ws.Range("C5").Value = "0603"
$ws.Range("D5").Value = "2.2µF ±10% 10V Ceramic Capacitor X5R 0603 (1608 Metric)"
$ws.Range("E5").Value = "SMD"
$ws.Range("F5").Value = "0603"
$ws.Range("G5").Value = "Samsung Electro-Mechanics"
$ws.Range("H5").Value = "CL10A225KP8NNNC"
$ws.Range("I5").Value = "Digi-Key"
$ws.Range("J5").Value = "1276-1085-1-ND"
$ws.Range("K5").Value = 0.1
$ws.Range("L5").Value = 1
$ws.Range("M5").Formula = "=K5*L5"
$ws.Range("N5").Value = 0.034
$ws.Range("O5").Formula = "=N5*L5"

# ---- row 6: 0.05" spaced header pins ----
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "0.05' spaced Header Pins"
$ws.Range("C6").Value = "J1"
$ws.Range("D6").Value = "Connector Header Through Hole 50 position 0.050"" (1.27mm)"
$ws.Range("E6").Value = "THT"
$ws.Range("F6").Value = "NA"
$ws.Range("G6").Value = "Sullins Connector Solutions"
$ws.Range("H6").Value = "GRPB501VWVN-RC"
$ws.Range("I6").Value = "Digi-Key"
$ws.Range("J6").Value = "S9014E-50-ND"
$ws.Range("K6").Value = 3.11
$ws.Range("L6").Value = 1
$ws.Range("M6").Formula = "=K6*L6"
$ws.Range("N6").Value = 2.4912
$ws.Range("O6").Formula = "=N6*L6"
